$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1212.7715
$ws.Range("I8").Value = 139.625
$ws.Range("J8").Value = 1530.7407
$ws.Range("K8").Value = 418.875
$ws.Range("L8").Value = 4592.2221
$ws.Range("M8").Value = -279.875
$ws.Range("N8").Value = -4870.2221

$ws.Range("H53").Value = 285.04544
$ws.Range("I53").Value = 224.15384
$ws.Range("J53").Value = 373
$ws.Range("K53").Value = 224.15384
$ws.Range("L53").Value = 373
$ws.Range("M53").Value = 412.84616
$ws.Range("N53").Value = -1647

$ws.Range("H74").Value = 3842.5
$ws.Range("I74").Value = 3766.6667
$ws.Range("J74").Value = 3888
$ws.Range("K74").Value = 3766.6667
$ws.Range("L74").Value = 3888
$ws.Range("M74").Value = -2830.6667
$ws.Range("N74").Value = -5760

$ws.Range("H77").Value = 3842.5
$ws.Range("I77").Value = 3766.6667
$ws.Range("J77").Value = 3888
$ws.Range("K77").Value = 18833.3335
$ws.Range("L77").Value = 19440
$ws.Range("M77").Value = -14153.3335
$ws.Range("N77").Value = -28800

$ws.Range("H92").Value = 350.5
$ws.Range("I92").Value = 352.86206
$ws.Range("J92").Value = 342.8889
$ws.Range("K92").Value = 352.86206
$ws.Range("L92").Value = 342.8889
$ws.Range("M92").Value = 895.1379400000001
$ws.Range("N92").Value = -2838.8889

$ws.Range("H112").Value = 1604.1177
$ws.Range("I112").Value = 2103.75
$ws.Range("J112").Value = 1160
$ws.Range("K112").Value = 6311.25
$ws.Range("L112").Value = 3480
$ws.Range("M112").Value = -5203.25
$ws.Range("N112").Value = -5696

$ws.Range("H115").Value = 1280.6875
$ws.Range("I115").Value = 490.0909
$ws.Range("J115").Value = 3020
$ws.Range("K115").Value = 1470.2727
$ws.Range("L115").Value = 9060
$ws.Range("M115").Value = 96.72730000000001
$ws.Range("N115").Value = -12194

$ws.Range("H118").Value = 2434.1738
$ws.Range("I118").Value = 487.8889
$ws.Range("J118").Value = 3685.3572
$ws.Range("K118").Value = 1463.6667
$ws.Range("L118").Value = 11056.0716
$ws.Range("M118").Value = 193.3333
$ws.Range("N118").Value = -14370.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27779774
$ws.Range("I2").Value = 41667412
$ws.Range("J2").Value = 4500
$ws.Range("K2").Value = 41667412
$ws.Range("L2").Value = 4500
$ws.Range("M2").Value = -41667299
$ws.Range("N2").Value = -4726

$ws.Range("H32").Value = 18595.982
$ws.Range("I32").Value = 15219.341
$ws.Range("J32").Value = 29208.285
$ws.Range("K32").Value = 15219.341
$ws.Range("L32").Value = 29208.285
$ws.Range("M32").Value = -14932.341
$ws.Range("N32").Value = -29782.285

$ws.Range("H45").Value = 1789.7059
$ws.Range("I45").Value = 1160.7333
$ws.Range("J45").Value = 6507
$ws.Range("K45").Value = 1160.7333
$ws.Range("L45").Value = 6507
$ws.Range("M45").Value = -783.7333000000001
$ws.Range("N45").Value = -7261

$ws.Range("H102").Value = 3246.35
$ws.Range("I102").Value = 2120.375
$ws.Range("J102").Value = 7750.25
$ws.Range("K102").Value = 2120.375
$ws.Range("L102").Value = 7750.25
$ws.Range("M102").Value = -498.375
$ws.Range("N102").Value = -10994.25

$ws.Range("H116").Value = 27779774
$ws.Range("I116").Value = 41667412
$ws.Range("J116").Value = 4500
$ws.Range("K116").Value = 41667412
$ws.Range("L116").Value = 4500
$ws.Range("M116").Value = -41665118
$ws.Range("N116").Value = -9088

$ws.Range("H132").Value = 2706.5757
$ws.Range("I132").Value = 2013.3636
$ws.Range("J132").Value = 4093
$ws.Range("K132").Value = 6040.0908
$ws.Range("L132").Value = 12279
$ws.Range("M132").Value = -3510.0908
$ws.Range("N132").Value = -17339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27779774
$ws.Range("I3").Value = 41667412
$ws.Range("J3").Value = 4500
$ws.Range("K3").Value = 41667412
$ws.Range("L3").Value = 4500
$ws.Range("M3").Value = -41667298
$ws.Range("N3").Value = -4728

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H94").Value = 663
$ws.Range("I94").Value = 535.48
$ws.Range("J94").Value = 1725.6666
$ws.Range("K94").Value = 535.48
$ws.Range("L94").Value = 1725.6666
$ws.Range("M94").Value = -84.48000000000002
$ws.Range("N94").Value = -2627.6666

$ws.Range("H96").Value = 21474.75
$ws.Range("I96").Value = 8949.5
$ws.Range("J96").Value = 34000
$ws.Range("K96").Value = 8949.5
$ws.Range("L96").Value = 34000
$ws.Range("M96").Value = -6203.5
$ws.Range("N96").Value = -39492

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

$ws.Range("H105").Value = 2069.25
$ws.Range("I105").Value = 1826.4117
$ws.Range("K105").Value = 1826.4117
$ws.Range("M105").Value = -79.41170000000011

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4403.0557
$ws.Range("I94").Value = 8284.799999999999
$ws.Range("J94").Value = 2910.077
$ws.Range("K94").Value = 8284.799999999999
$ws.Range("L94").Value = 2910.077
$ws.Range("M94").Value = -7833.799999999999
$ws.Range("N94").Value = -3812.077

$ws.Range("H107").Value = 1248.1389
$ws.Range("J107").Value = 1459.3636
$ws.Range("L107").Value = 1459.3636
$ws.Range("N107").Value = -5299.3636

$ws.Range("H134").Value = 2807.7715
$ws.Range("I134").Value = 2019.6666
$ws.Range("J134").Value = 4527.273
$ws.Range("K134").Value = 6058.9998
$ws.Range("L134").Value = 13581.819
$ws.Range("M134").Value = -3523.9998
$ws.Range("N134").Value = -18651.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3432
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 3702.2222
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 11106.6666
$ws.Range("M94").Value = -2324
$ws.Range("N94").Value = -12458.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3560
$ws.Range("I80").Value = 3616.6667
$ws.Range("J80").Value = 3220
$ws.Range("K80").Value = 3616.6667
$ws.Range("L80").Value = 3220
$ws.Range("M80").Value = -2618.6667
$ws.Range("N80").Value = -5216

$ws.Range("H83").Value = 3560
$ws.Range("I83").Value = 3616.6667
$ws.Range("J83").Value = 3220
$ws.Range("K83").Value = 18083.3335
$ws.Range("L83").Value = 16100
$ws.Range("M83").Value = -13091.3335
$ws.Range("N83").Value = -26084

$ws.Range("H132").Value = 4696.853
$ws.Range("I132").Value = 5624.1665
$ws.Range("J132").Value = 3653.625
$ws.Range("K132").Value = 16872.4995
$ws.Range("L132").Value = 10960.875
$ws.Range("M132").Value = -14342.4995
$ws.Range("N132").Value = -16020.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1219.3549
$ws.Range("I46").Value = 892.3077
$ws.Range("K46").Value = 892.3077
$ws.Range("M46").Value = -704.3077

$ws.Range("H106").Value = 26636.7
$ws.Range("J106").Value = 26636.7
$ws.Range("L106").Value = 26636.7
$ws.Range("N106").Value = -29160.7

$ws.Range("H132").Value = 2529.551
$ws.Range("I132").Value = 1756
$ws.Range("J132").Value = 4125
$ws.Range("K132").Value = 5268
$ws.Range("L132").Value = 12375
$ws.Range("M132").Value = -2738
$ws.Range("N132").Value = -17435

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H132").Value = 2186476
$ws.Range("I132").Value = 2635336.5
$ws.Range("J132").Value = 54388.625
$ws.Range("K132").Value = 7906009.5
$ws.Range("L132").Value = 163165.875
$ws.Range("M132").Value = -7903479.5
$ws.Range("N132").Value = -168225.875

